$d = $word.ActiveDocument

# Locate the "Ben Bar" run sequence (name title on the CV header table).
$rng = $d.Content.Duplicate
$rng.Find.ClearFormatting()
$found = $rng.Find.Execute("Ben Bar", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

# Work on the whole containing paragraph so we can rewrite its OOXML with
# InsertXML (which replaces the full paragraph once the range touches it),
# while keeping the original <w:pPr> and the first three runs byte-for-byte
# identical and only appending a new trailing run "rrr".
$para = $rng.Paragraphs(1).Range

$newParaXml = '<w:p w14:paraId="7E3D6BC3" w14:textId="733FAC1E" w:rsidR="00BF2955" w:rsidRPr="003A5831" w:rsidRDefault="00BF2955" w:rsidP="003462C6">' + `
  '<w:pPr><w:pStyle w:val="1"/><w:jc w:val="left"/><w:rPr><w:rFonts w:ascii="Cambria" w:hAnsi="Cambria" w:cstheme="majorBidi"/><w:b w:val="0"/><w:bCs/><w:caps w:val="0"/><w:color w:val="A6A6A6" w:themeColor="background1" w:themeShade="A6"/><w:sz w:val="52"/><w:szCs w:val="52"/><w:lang w:bidi="he-IL"/></w:rPr></w:pPr>' + `
  '<w:r w:rsidRPr="003A5831"><w:rPr><w:rFonts w:ascii="Cambria" w:hAnsi="Cambria" w:cs="Poppins"/><w:b w:val="0"/><w:bCs/><w:caps w:val="0"/><w:color w:val="FFFFFF" w:themeColor="background1"/><w:sz w:val="52"/><w:szCs w:val="52"/></w:rPr><w:t xml:space="preserve">    </w:t></w:r>' + `
  '<w:r w:rsidRPr="003A5831"><w:rPr><w:rFonts w:ascii="Cambria" w:hAnsi="Cambria" w:cstheme="majorBidi"/><w:b w:val="0"/><w:bCs/><w:caps w:val="0"/><w:color w:val="FFFFFF" w:themeColor="background1"/><w:sz w:val="52"/><w:szCs w:val="52"/></w:rPr><w:t>Ben Ba</w:t></w:r>' + `
  '<w:r w:rsidR="00FB7AED" w:rsidRPr="003A5831"><w:rPr><w:rFonts w:ascii="Cambria" w:hAnsi="Cambria" w:cstheme="majorBidi"/><w:b w:val="0"/><w:bCs/><w:caps w:val="0"/><w:color w:val="FFFFFF" w:themeColor="background1"/><w:sz w:val="52"/><w:szCs w:val="52"/></w:rPr><w:t>r</w:t></w:r>' + `
  '<w:r><w:rPr><w:rFonts w:ascii="Cambria" w:hAnsi="Cambria" w:cstheme="majorBidi"/><w:b w:val="0"/><w:bCs/><w:caps w:val="0"/><w:color w:val="FFFFFF" w:themeColor="background1"/><w:sz w:val="52"/><w:szCs w:val="52"/></w:rPr><w:t>rrr</w:t></w:r>' + `
  '</w:p>'

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
      '<pkg:xmlData>' + `
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
          '<w:body>' + $newParaXml + '</w:body>' + `
        '</w:document>' + `
      '</pkg:xmlData>' + `
    '</pkg:part>' + `
  '</pkg:package>'

$para.InsertXML($xml)
